$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.473.58'
$ws.Range("D3").Value = '1.913.06'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9986'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9988'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  +2.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2893'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06727'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '111.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.12'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.12%  '
$ws.Range("D12").Value = '1.908.17'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07558'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.274'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6723'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '287.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = '30.476.57'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007603'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9990'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("D21").Value = '2.164.07'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.475'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9984'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.407'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.473'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1055'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.402'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.182'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.052'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04986'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7301'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.136'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9988'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.719'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02036'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.664'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '110.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.015'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4437'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8677'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.834'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9988'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  +0.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.366'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.261'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1241'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.36%  '
